$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1, styled like the other header cells (copy E1's style)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Data rows with time_taken timestamps (plain text, unstyled like B2:E11)
$ws.Range("F2").Value = "2021-10-05 13:38:31.625977"
$ws.Range("F3").Value = "2021-10-05 13:38:31.625987"
$ws.Range("F4").Value = "2021-10-05 13:38:31.625991"
$ws.Range("F5").Value = "2021-10-05 13:38:31.625993"
$ws.Range("F6").Value = "2021-10-05 13:38:31.625996"
$ws.Range("F7").Value = "2021-10-05 13:38:31.625998"
$ws.Range("F8").Value = "2021-10-05 13:38:31.626001"
$ws.Range("F9").Value = "2021-10-05 13:38:31.626003"
$ws.Range("F10").Value = "2021-10-05 13:38:31.626006"
$ws.Range("F11").Value = "2021-10-05 13:38:31.626009"
